# Trade #13 closed at 2026-02-17 07:58:49 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 13
$wsSummary.Range("B9").Value = 38.46

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 13
$wsStatus.Range("G4").Value = 38.46

# --- New trade row data (trade #13) ---
$tradeNum = 13
$tradeDate = "2026-02-17"
$tradeTime = "07:58:42"
$tradeStrategy = "MarketMaking"
$tradeSide = "DOWN"
$entryPrice = 0.95
$exitPrice = 0.95
$tradeStatus = "CLOSED"
$pnlPct = 0
$pnlDollar = 0
$capitalAfter = 99.98
$entrySlippage = 0
$exitSlippage = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.13

function Add-TradeRow($ws) {
    $ws.Cells.Item(14, 1).Value = $tradeNum
    $ws.Cells.Item(14, 2).Value = "'" + $tradeDate
    $ws.Cells.Item(14, 3).Value = $tradeTime
    $ws.Cells.Item(14, 4).Value = $tradeStrategy
    $ws.Cells.Item(14, 5).Value = $tradeSide
    $ws.Cells.Item(14, 6).Value = $entryPrice
    $ws.Cells.Item(14, 7).Value = $exitPrice
    $ws.Cells.Item(14, 8).Value = $tradeStatus
    $ws.Cells.Item(14, 9).Value = $pnlPct
    $ws.Cells.Item(14, 10).Value = $pnlDollar
    $ws.Cells.Item(14, 11).Value = $capitalAfter
    $ws.Cells.Item(14, 12).Value = $entrySlippage
    $ws.Cells.Item(14, 13).Value = $exitSlippage
    $ws.Cells.Item(14, 14).Value = $confidence
    $ws.Cells.Item(14, 15).Value = $entryReason
    $ws.Cells.Item(14, 16).Value = $exitReason
    $ws.Cells.Item(14, 17).Value = $duration
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
